{"js": "// Update the date label and the 25 two-digit multiplication problems to\n// the new values from the commit (\"output generated at c8c62b6\").\nconst replacements = [\n  [\"2025-05-07 Wednesday\", \"2025-05-08 Thursday\"],\n  [\"67\u00d797=\", \"99\u00d790=\"],\n  [\"42\u00d717=\", \"25\u00d713=\"],\n  [\"66\u00d739=\", \"13\u00d743=\"],\n  [\"17\u00d755=\", \"80\u00d713=\"],\n  [\"22\u00d713=\", \"81\u00d779=\"],\n  [\"88\u00d746=\", \"56\u00d717=\"],\n  [\"29\u00d781=\", \"22\u00d794=\"],\n  [\"97\u00d782=\", \"95\u00d725=\"],\n  [\"99\u00d723=\", \"93\u00d727=\"],\n  [\"37\u00d774=\", \"99\u00d738=\"],\n  [\"51\u00d739=\", \"60\u00d753=\"],\n  [\"98\u00d759=\", \"30\u00d759=\"],\n  [\"41\u00d770=\", \"39\u00d716=\"],\n  [\"24\u00d762=\", \"64\u00d730=\"],\n  [\"33\u00d765=\", \"71\u00d743=\"],\n  [\"12\u00d768=\", \"20\u00d716=\"],\n  [\"99\u00d733=\", \"48\u00d715=\"],\n  [\"23\u00d769=\", \"18\u00d712=\"],\n  [\"16\u00d773=\", \"32\u00d741=\"],\n  [\"52\u00d796=\", \"16\u00d780=\"],\n  [\"63\u00d799=\", \"68\u00d714=\"],\n  [\"67\u00d796=\", \"89\u00d716=\"],\n  [\"31\u00d775=\", \"46\u00d778=\"],\n  [\"51\u00d729=\", \"37\u00d754=\"],\n  [\"56\u00d753=\", \"70\u00d755=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and the 25 two-digit multiplication problems to\n# the new values from the commit (\"output generated at c8c62b6\").\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-07 Wednesday\", \"2025-05-08 Thursday\"),\n    @(\"67\u00d797=\", \"99\u00d790=\"),\n    @(\"42\u00d717=\", \"25\u00d713=\"),\n    @(\"66\u00d739=\", \"13\u00d743=\"),\n    @(\"17\u00d755=\", \"80\u00d713=\"),\n    @(\"22\u00d713=\", \"81\u00d779=\"),\n    @(\"88\u00d746=\", \"56\u00d717=\"),\n    @(\"29\u00d781=\", \"22\u00d794=\"),\n    @(\"97\u00d782=\", \"95\u00d725=\"),\n    @(\"99\u00d723=\", \"93\u00d727=\"),\n    @(\"37\u00d774=\", \"99\u00d738=\"),\n    @(\"51\u00d739=\", \"60\u00d753=\"),\n    @(\"98\u00d759=\", \"30\u00d759=\"),\n    @(\"41\u00d770=\", \"39\u00d716=\"),\n    @(\"24\u00d762=\", \"64\u00d730=\"),\n    @(\"33\u00d765=\", \"71\u00d743=\"),\n    @(\"12\u00d768=\", \"20\u00d716=\"),\n    @(\"99\u00d733=\", \"48\u00d715=\"),\n    @(\"23\u00d769=\", \"18\u00d712=\"),\n    @(\"16\u00d773=\", \"32\u00d741=\"),\n    @(\"52\u00d796=\", \"16\u00d780=\"),\n    @(\"63\u00d799=\", \"68\u00d714=\"),\n    @(\"67\u00d796=\", \"89\u00d716=\"),\n    @(\"31\u00d775=\", \"46\u00d778=\"),\n    @(\"51\u00d729=\", \"37\u00d754=\"),\n    @(\"56\u00d753=\", \"70\u00d755=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
